# Refresh the cryptos list: update the Price (D) and Volume(1h) (E) columns
# for rows 2-51. A leading apostrophe is used for Price values that would
# otherwise be auto-parsed as numbers (e.g. "561.79"), so they keep being
# stored as text just like in the original workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.101.55'
$ws.Range('E2').Value = '  -0.82%  '
$ws.Range('D3').Value = '3.150.54'
$ws.Range('E3').Value = '  +1.90%  '
$ws.Range('E4').Value = '  +0.19%  '
$ws.Range('D5').Value = "'561.79"
$ws.Range('E5').Value = '  +0.84%  '
$ws.Range('D6').Value = "'140.90"
$ws.Range('E6').Value = '  +0.07%  '
$ws.Range('E7').Value = '  +0.28%  '
$ws.Range('D8').Value = '3.136.69'
$ws.Range('E8').Value = '  +1.63%  '
$ws.Range('D9').Value = "'0.492"
$ws.Range('E9').Value = '  -0.36%  '
$ws.Range('D10').Value = "'6.69"
$ws.Range('E10').Value = '  +1.58%  '
$ws.Range('D11').Value = "'0.153"
$ws.Range('E11').Value = '  -2.17%  '
$ws.Range('E12').Value = '  -0.07%  '
$ws.Range('D13').Value = "'36.27"
$ws.Range('E13').Value = '  +0.50%  '
$ws.Range('E14').Value = '  -1.46%  '
$ws.Range('D15').Value = '3.658.17'
$ws.Range('E15').Value = '  +2.10%  '
$ws.Range('D16').Value = '64.256.56'
$ws.Range('E16').Value = '  -0.58%  '
$ws.Range('D17').Value = '3.151.95'
$ws.Range('E17').Value = '  +2.19%  '
$ws.Range('E18').Value = '  +0.56%  '
$ws.Range('D19').Value = "'509.93"
$ws.Range('E19').Value = '  +2.89%  '
$ws.Range('D20').Value = "'6.78"
$ws.Range('E20').Value = '  +0.99%  '
$ws.Range('D21').Value = "'13.93"
$ws.Range('E21').Value = '  +1.05%  '
$ws.Range('E22').Value = '  +2.86%  '
$ws.Range('D23').Value = "'7.40"
$ws.Range('E23').Value = '  +1.33%  '
$ws.Range('D24').Value = "'12.72"
$ws.Range('E24').Value = '  +0.93%  '
$ws.Range('D25').Value = "'78.43"
$ws.Range('E25').Value = '  -1.07%  '
$ws.Range('E26').Value = '  -0.03%  '
$ws.Range('D27').Value = "'8.70"
$ws.Range('E27').Value = '  +8.86%  '
$ws.Range('E28').Value = '  +2.82%  '
$ws.Range('E29').Value = '  -0.43%  '
$ws.Range('E30').Value = '  +0.12%  '
$ws.Range('D31').Value = "'26.47"
$ws.Range('E31').Value = '  +0.86%  '
$ws.Range('E32').Value = '  -2.75%  '
$ws.Range('D33').Value = "'1.13"
$ws.Range('E33').Value = '  -0.09%  '
$ws.Range('E34').Value = '  -5.22%  '
$ws.Range('E35').Value = '  -0.40%  '
$ws.Range('D36').Value = "'5.32"
$ws.Range('E36').Value = '  -3.33%  '
$ws.Range('D37').Value = "'53.65"
$ws.Range('E37').Value = '  +1.70%  '
$ws.Range('D38').Value = "'0.0425"
$ws.Range('E38').Value = '  +3.69%  '
$ws.Range('D39').Value = '3.155.92'
$ws.Range('E39').Value = '  +6.07%  '
$ws.Range('D40').Value = "'0.0814"
$ws.Range('E40').Value = '  +1.65%  '
$ws.Range('E41').Value = '  +1.81%  '
$ws.Range('D42').Value = "'2.73"
$ws.Range('E42').Value = '  -7.54%  '
$ws.Range('D43').Value = "'8.21"
$ws.Range('E43').Value = '  -1.57%  '
$ws.Range('E44').Value = '  +6.12%  '
$ws.Range('D45').Value = "'2.17"
$ws.Range('E45').Value = '  +2.08%  '
$ws.Range('E46').Value = '  -0.04%  '
$ws.Range('D47').Value = "'122.11"
$ws.Range('E47').Value = '  +1.57%  '
$ws.Range('D48').Value = "'24.83"
$ws.Range('E48').Value = '  -1.69%  '
$ws.Range('E49').Value = '  -1.55%  '
$ws.Range('D50').Value = '0.0₃0512'
$ws.Range('E50').Value = '  -6.53%  '
$ws.Range('E51').Value = '  -1.23%  '
